# Refresh cryptocurrency price/volume data (GitHub Actions scheduled update)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "58.859.53"
$ws.Range("E2").Value = "  +2.98%  "

# Row 3
$ws.Range("D3").Value = "2.589.61"
$ws.Range("E3").Value = "  +1.57%  "

# Row 4
$ws.Range("E4").Value = "  +0.06%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "520.69"
$ws.Range("E5").Value = "  +0.86%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "140.11"
$ws.Range("E6").Value = "  -0.67%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.998"
$ws.Range("E7").Value = "  -0.04%  "

# Row 8
$ws.Range("E8").Value = "  +1.33%  "

# Row 9
$ws.Range("D9").Value = "2.602.56"
$ws.Range("E9").Value = "  +1.55%  "

# Row 10
$ws.Range("E10").Value = "  -1.36%  "

# Row 11
$ws.Range("E11").Value = "  +1.79%  "

# Row 12
$ws.Range("E12").Value = "  +2.92%  "

# Row 13
$ws.Range("E13").Value = "  +2.55%  "

# Row 14
$ws.Range("D14").Value = "3.046.66"
$ws.Range("E14").Value = "  +1.59%  "

# Row 15
$ws.Range("D15").Value = "58.824.27"
$ws.Range("E15").Value = "  +2.90%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "20.48"
$ws.Range("E16").Value = "  +2.43%  "

# Row 17
$ws.Range("D17").Value = "2.608.08"
$ws.Range("E17").Value = "  +1.44%  "

# Row 18
$ws.Range("E18").Value = "  +0.58%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "339.27"
$ws.Range("E19").Value = "  +2.50%  "

# Row 21
$ws.Range("E21").Value = "  +1.35%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.53"
$ws.Range("E22").Value = "  +6.14%  "

# Row 23
$ws.Range("E23").Value = "  +0.21%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "66.30"
$ws.Range("E24").Value = "  +2.33%  "

# Row 25
$ws.Range("E25").Value = "  -0.22%  "

# Row 26
$ws.Range("E26").Value = "  +1.73%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.998"
$ws.Range("E27").Value = "  +0.04%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.09"
$ws.Range("E28").Value = "  +2.81%  "

# Row 29
$ws.Range("E29").Value = "  +0.04%  "

# Row 30
$ws.Range("D30").Value = "0.0₃0725"
$ws.Range("E30").Value = "  -1.39%  "

# Row 31
$ws.Range("E31").Value = "  -4.67%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "18.79"
$ws.Range("E32").Value = "  +1.80%  "

# Row 33
$ws.Range("E33").Value = "  -0.01%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "148.69"
$ws.Range("E34").Value = "  +0.19%  "

# Row 35
$ws.Range("E35").Value = "  +0.78%  "

# Row 36
$ws.Range("E36").Value = "  -0.10%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "36.26"
$ws.Range("E37").Value = "  +1.88%  "

# Row 38
$ws.Range("E38").Value = "  +1.96%  "

# Row 39
$ws.Range("E39").Value = "  +1.98%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.823"
$ws.Range("E40").Value = "  -1.17%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.51"
$ws.Range("E41").Value = "  +1.45%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.998"
$ws.Range("E42").Value = "  -0.13%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "276.37"
$ws.Range("E43").Value = "  +4.22%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "10.73"
$ws.Range("E44").Value = "  +0.90%  "

# Row 45
$ws.Range("E45").Value = "  +0.24%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.587"
$ws.Range("E46").Value = "  +1.41%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0524"
$ws.Range("E47").Value = "  +1.91%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "18.64"
$ws.Range("E48").Value = "  -0.10%  "

# Row 49
$ws.Range("D49").Value = "1.986.70"
$ws.Range("E49").Value = "  +1.36%  "

# Row 50
$ws.Range("E50").Value = "  +1.19%  "

# Row 51
$ws.Range("E51").Value = "  -0.49%  "
